$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.935.20"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").Value = "2.355.54"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +4.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.52"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +11.21%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +19.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.101"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "29.10"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.47%  "
$ws.Range("E12").Value = "  +2.73%  "
$ws.Range("D13").Value = "2.706.42"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.902"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.57%  "
$ws.Range("D17").Value = "2.377.21"
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("D18").Value = "43.907.06"
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("E19").Value = "  +5.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "77.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "254.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.86%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E25").Value = "  +4.15%  "
$ws.Range("E26").Value = "  +6.46%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  +1.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("E30").Value = "  +7.40%  "
$ws.Range("E31").Value = "  +1.95%  "
$ws.Range("E32").Value = "  +5.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0719"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.17%  "
$ws.Range("E37").Value = "  -1.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("E39").Value = "  +7.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.91%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("E43").Value = "  +4.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0982"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.36%  "
$ws.Range("E45").Value = "  +1.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("E48").Value = "  +12.38%  "
$ws.Range("E49").Value = "  +4.10%  "
$ws.Range("D50").Value = "1.436.70"
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("E51").Value = "  +1.29%  "
